$d = $word.ActiveDocument

# 1. Merge "...as soon as p" + "ossible." runs into a single run's text.
$d.Content.Find.Execute("as soon as possible.", $true, $false, $false, $false, $false, $true, 1, $false, "as soon as possible.", 2)

# 2. Rewrite the "Automated Testing Instructions" paragraph, splitting it
#    into the specific sequence of runs shown by the diff. The replacement
#    text is spliced in as one run (inheriting the existing run formatting),
#    then we "cut" it into the target run boundaries by toggling Bold on/off
#    (restoring the original value) at each boundary -- toggling a formatting
#    property on a sub-range forces the engine to materialize a run split
#    there without altering the visible formatting.
$full = $d.Content
$oldText = "We have done a lot of manual testing on the Emulator. We have planned to develop and release automated testing by Release 2. Therefore, following our initial plan, automated testing is not implemented in this release."
$full.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $full.Start

$newText = "We did a lot of manual testing on the Emulator. Everything is working so far, and we have planned to develop and release an automated testing system by Release 2. Therefore, following our initial plan, automated testing is not implemented in this release."
$full.Text = $newText
$newEnd = $start + $newText.Length

# Offsets (relative to $start) of the boundaries between the 10 target runs:
# "W" | "e " | "did " | "a lot of manual testing on the Emulator. " |
# "Everything is working so far, and w" | "e have planned to develop and release " |
# "an " | "automated testing " | "system by " |
# "Release 2. Therefore, following our initial plan, automated testing is not implemented in this release."
$boundaries = @(1, 3, 7, 48, 83, 121, 124, 142, 152) | Sort-Object -Descending
foreach ($b in $boundaries) {
    $p = $start + $b
    $rb = $d.Range($p, $newEnd)
    $rb.Bold = 1
    $rb.Bold = 0
}
